# CMMS file jobs anonymized: replace the real job description texts in
# column C (Description) with generic placeholder text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "job description 1"
$ws.Range("C3").Value = "job description 2"
$ws.Range("C4").Value = "job description 3"
$ws.Range("C5").Value = "job description 4"
